# TMTT0017339_EventExpense_VerifyTheCFExpenseRequestFunctionalityAsApprover - 21 Dec 2023
# Updates the stored Outlook test credentials (username/password) and refreshes
# the related sheet formatting (selection, column widths, hyperlink) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update credential values -------------------------------------------------
$ws.Range("A2").Value = "Sahil.Mittal0207@hl.com"
$ws.Range("B2").Value = "Yankee@123456"

# --- Refresh the hyperlinks (the underlying mailto targets are untouched by
# the authored change - only the stale "display text" override on B2 goes
# away). Re-add both so the collection keeps exactly two entries tied to the
# same rId order as before.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:vkumar0427@hl.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Yankee@12345")

# Restore the "Hyperlink" cell style that Add() re-stamps with a fresh (but
# equivalent) style record, so the underlying style indexes stay untouched.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"

# --- Column widths (auto-fit grew/shrank to match the new text lengths) ------
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 13.76

# --- Active selection moved from C2 to B2 -------------------------------------
$ws.Range("B2").Select() | Out-Null
